$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Zener diode line item (D4, BZT52C2V0S-7-F, Diodes Inc.) occupied row 9
# of the BOM. Delete it entirely; every row beneath it shifts up by one.
$ws.Rows("9").Delete()

# The workbook-level defined name "H0FR70" spans the whole BOM table and
# must shrink along with the removed row.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!H0FR70") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$H`$22"
    }
}

# After deleting the row, Excel leaves the new row 9 (which used to be row 10)
# selected as a whole row.
$ws.Range("A9:XFD9").Select()
